# Populate the GO node IDs in column A (replacing the placeholder
# "<built-in function id>" text) and fix up the childnodes counts
# in column D for rows 4 and 33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ids = @{
    2  = "GO:0000045"
    3  = "GO:0000421"
    4  = "GO:0016236"
    5  = "GO:0016237"
    6  = "GO:0016240"
    7  = "GO:0016243"
    8  = "GO:0030399"
    9  = "GO:0032258"
    10 = "GO:0034423"
    11 = "GO:0044753"
    12 = "GO:0044754"
    13 = "GO:0045771"
    14 = "GO:0045772"
    15 = "GO:0048102"
    16 = "GO:0061709"
    17 = "GO:0061739"
    18 = "GO:0061753"
    19 = "GO:0061906"
    20 = "GO:0061908"
    21 = "GO:0061909"
    22 = "GO:0061910"
    23 = "GO:0097635"
    24 = "GO:0097636"
    25 = "GO:0097637"
    26 = "GO:0098792"
    27 = "GO:0120095"
    28 = "GO:1901096"
    29 = "GO:1901097"
    30 = "GO:1901098"
    31 = "GO:1901245"
    32 = "GO:1902902"
    33 = "GO:1905037"
    34 = "GO:1990316"
    35 = "GO:1990462"
    36 = "GO:2000785"
}

foreach ($row in $ids.Keys) {
    $ws.Cells.Item($row, 1).Value = $ids[$row]
}

# Update the childnodes counts that changed alongside the ids.
$ws.Cells.Item(4, 4).Value = 10
$ws.Cells.Item(33, 4).Value = 5
